$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.545.88'
$ws.Range('E2').Value = '  -1.10%  '

$ws.Range('D3').Value = '3.542.78'
$ws.Range('E3').Value = '  -1.87%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.18'
$ws.Range('E5').Value = '  -3.31%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.64'
$ws.Range('E6').Value = '  -0.32%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.610'

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  -0.53%  '

$ws.Range('E10').Value = '  -2.81%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.54'
$ws.Range('E11').Value = '  -2.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000288'
$ws.Range('E12').Value = '  -5.45%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.22'
$ws.Range('E13').Value = '  -3.64%  '

$ws.Range('D14').Value = '4.099.89'
$ws.Range('E14').Value = '  -1.93%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '666.29'
$ws.Range('E15').Value = '  +11.75%  '

$ws.Range('D16').Value = '69.535.57'
$ws.Range('E16').Value = '  -1.35%  '

$ws.Range('D17').Value = '3.533.06'
$ws.Range('E17').Value = '  -1.88%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.48'
$ws.Range('E18').Value = '  -4.35%  '

$ws.Range('E19').Value = '  -0.95%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.40'
$ws.Range('E20').Value = '  -3.34%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.962'
$ws.Range('E21').Value = '  -3.62%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.03'
$ws.Range('E22').Value = '  +1.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.40'
$ws.Range('E23').Value = '  +4.41%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.27'
$ws.Range('E24').Value = '  +1.98%  '

$ws.Range('E25').Value = '  -5.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.93'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.10'
$ws.Range('E27').Value = '  -6.35%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.55'
$ws.Range('E28').Value = '  -0.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.21'
$ws.Range('E29').Value = '  -2.03%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.35'
$ws.Range('E30').Value = '  -8.11%  '

$ws.Range('E31').Value = '  -6.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.75'
$ws.Range('E32').Value = '  -4.52%  '

$ws.Range('E33').Value = '  -5.62%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '61.67'
$ws.Range('E34').Value = '  -2.41%  '

$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.81'
$ws.Range('E35').Value = '  +7.92%  '

$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.767.14'
$ws.Range('E36').Value = '  -3.24%  '

$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '0.0₃0811'
$ws.Range('E37').Value = '  -10.22%  '

$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '498.52'
$ws.Range('E39').Value = '  -7.38%  '

$ws.Range('E40').Value = '  -7.83%  '

$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.370'
$ws.Range('E41').Value = '  -5.63%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  +0.57%  '

$ws.Range('E43').Value = '  -6.56%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0448'
$ws.Range('E44').Value = '  -1.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  -0.31%  '

$ws.Range('E46').Value = '  -0.08%  '

$ws.Range('E47').Value = '  -3.04%  '

$ws.Range('E48').Value = '  -0.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.32'
$ws.Range('E49').Value = '  -3.50%  '

$ws.Range('E50').Value = '  +19.32%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.70'
$ws.Range('E51').Value = '  +62.59%  '

